$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '43.955.43'
$ws.Range('E2').Value = '  +2.56%  '
$ws.Range('D3').Value = '2.342.57'
$ws.Range('E3').Value = '  +2.47%  '
$ws.Range('E4').Value = '  +0.35%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '313.45'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  -0.58%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '108.93'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  +4.19%  '
$ws.Range('E7').Value = '  +1.54%  '
$ws.Range('E8').Value = '  -0.14%  '
$ws.Range('E9').Value = '  +2.73%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '41.16'
$ws.Range('D10').Style = "Normal"
$ws.Range('E10').Value = '  +4.16%  '
$ws.Range('E11').Value = '  +2.13%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '8.59'
$ws.Range('D12').Style = "Normal"
$ws.Range('E12').Value = '  +2.10%  '
$ws.Range('E13').Value = '  +1.45%  '
$ws.Range('E14').Value = '  -0.98%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '15.53'
$ws.Range('D15').Style = "Normal"
$ws.Range('E15').Value = '  +2.24%  '
$ws.Range('D16').Value = '2.697.04'
$ws.Range('E16').Value = '  +2.32%  '
$ws.Range('D17').Value = '2.328.71'
$ws.Range('E17').Value = '  +1.74%  '
$ws.Range('D18').Value = '43.838.41'
$ws.Range('E18').Value = '  +2.47%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '7.57'
$ws.Range('D19').Style = "Normal"
$ws.Range('E19').Value = '  +2.14%  '
$ws.Range('E20').Value = '  +1.89%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '13.04'
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value = '  -2.98%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '74.32'
$ws.Range('D22').Style = "Normal"
$ws.Range('E22').Value = '  +0.62%  '
$ws.Range('E23').Value = '  -1.27%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '269.60'
$ws.Range('D24').Style = "Normal"
$ws.Range('E24').Value = '  +2.83%  '
$ws.Range('E25').Value = '  +3.90%  '
$ws.Range('E26').Value = '  -0.31%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '7.63'
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value = '  +6.81%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '11.17'
$ws.Range('D28').Style = "Normal"
$ws.Range('E28').Value = '  +2.86%  '
$ws.Range('E29').Value = '  -2.39%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '38.89'
$ws.Range('D30').Style = "Normal"
$ws.Range('E30').Value = '  +4.43%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '22.70'
$ws.Range('D31').Style = "Normal"
$ws.Range('E31').Value = '  +2.00%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '168.39'
$ws.Range('D32').Style = "Normal"
$ws.Range('E32').Value = '  +1.05%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '0.0888'
$ws.Range('D33').Style = "Normal"
$ws.Range('E33').Value = '  +1.52%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '2.77'
$ws.Range('D34').Style = "Normal"
$ws.Range('E35').Value = '  +1.96%  '
$ws.Range('E36').Value = '  +5.15%  '
$ws.Range('E37').Value = '  -0.70%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.0367'
$ws.Range('D38').Style = "Normal"
$ws.Range('E38').Value = '  +5.12%  '
$ws.Range('E39').Value = '  +8.45%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '3.80'
$ws.Range('D40').Style = "Normal"
$ws.Range('E40').Value = '  -0.64%  '
$ws.Range('E41').Value = '  +8.68%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '104.36'
$ws.Range('D42').Style = "Normal"
$ws.Range('E42').Value = '  +12.88%  '
$ws.Range('E43').Value = '  +3.13%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '71.78'
$ws.Range('D44').Style = "Normal"
$ws.Range('E44').Value = '  +2.89%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '13.34'
$ws.Range('D45').Style = "Normal"
$ws.Range('E45').Value = '  +9.22%  '
$ws.Range('E46').Value = '  -0.14%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '114.13'
$ws.Range('D47').Style = "Normal"
$ws.Range('E47').Value = '  +0.58%  '
$ws.Range('D48').Value = '1.671.34'
$ws.Range('E48').Value = '  -3.36%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '76.95'
$ws.Range('D49').Style = "Normal"
$ws.Range('E49').Value = '  -2.71%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '8.99'
$ws.Range('D50').Style = "Normal"
$ws.Range('E50').Value = '  +2.69%  '
$ws.Range('B51').Value = 'TheGraph'
$ws.Range('C51').Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '0.216'
$ws.Range('D51').Style = "Normal"
$ws.Range('E51').Value = '  +13.84%  '
